# Insert a new data row at row 61 (pushing existing rows 61-125 down to
# 62-126) and populate the newly inserted row with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(61).Insert()

$ws.Cells.Item(61, 1).Value = 1
$ws.Cells.Item(61, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(61, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(61, 4).Value = 44966
$ws.Cells.Item(61, 5).Value = 15
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100102
$ws.Cells.Item(61, 8).Value = "Cítricos"
$ws.Cells.Item(61, 9).Value = 100102005
$ws.Cells.Item(61, 10).Value = "Naranja"
$ws.Cells.Item(61, 11).Value = "Valencia"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 300
$ws.Cells.Item(61, 14).Value = 900
$ws.Cells.Item(61, 15).Value = 950
$ws.Cells.Item(61, 16).Value = 925
$ws.Cells.Item(61, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 925
$ws.Cells.Item(61, 20).Value = 1
